$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "file_inspect_title"
$ws.Range("B24").Value = "File Inspect: {0}"

$ws.Range("A25").Select()
